# Remove sample codes from template
$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": clear the example QBIC sample id codes in A2:A6 ---
$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("A2").Value = ""
$metadata.Range("A3").Value = ""
$metadata.Range("A4").Value = ""
$metadata.Range("A5").Value = ""
$metadata.Range("A6").Value = ""

# --- Sheet "Property information": pluralize "QBIC sample id" -> "QBIC sample ids"
#     and fix the "measuremed" -> "measured" typo ---
$propInfo = $wb.Worksheets.Item("Property information")
$propInfo.Range("A2").Value = "QBIC sample ids"
$propInfo.Range("C11").Value = "In case of pooled sample get measured, indicate with a common sample group label for samples that are in the same measurement. Entries that share the same pool label will be combined as one measurement"

# --- Sheet "Allowed-Values": keep header consistent with the pluralized name ---
$allowed = $wb.Worksheets.Item("Allowed-Values")
$allowed.Range("A1").Value = "QBIC sample ids*"
